$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the now-superfluous last blank row so row count goes from 31 to 30.
$ws.Rows("31:31").Delete() | Out-Null

# Rewrite the header row (row 1) completely: idx / idx2 / Name / Date Start / Date End
# for the descriptive columns, then unit headers for the measurement columns.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

$ws.Range("F1:K1").Font.Size = 9
$ws.Range("F1:K1").Font.Name = "Arial"

# Move the two canton data rows up (they used to start at row 3, now start at row 2)
# since the old units row (row 2) is gone.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 107500
$ws.Range("C2").Value = "Kubel"
$ws.Range("D2").Value = 1900
$ws.Range("E2").Value = 1976
$ws.Range("F2").Value = 18
$ws.Range("G2").Value = 1.48
$ws.Range("H2").Value = 1.37
$ws.Range("I2").Value = 1.4
$ws.Range("J2").Value = 1.66
$ws.Range("K2").Value = 3.06

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 107400
$ws.Range("C3").Value = "Wasserauen"
$ws.Range("D3").Value = 1905
$ws.Range("E3").Value = 2005
$ws.Range("F3").Value = 1.3
$ws.Range("G3").Value = 2.5
$ws.Range("H3").Value = 2.5
$ws.Range("I3").Value = 1.48
$ws.Range("J3").Value = 6.22
$ws.Range("K3").Value = 7.7

# Row 4 used to hold the second canton's data; blank it out like the filler rows below.
$ws.Range("B4:K4").ClearContents() | Out-Null

$ws.Range("A2:K2").Select() | Out-Null
